# "Elimina EC anteriores y se agregan nuevos, se modifica base de datos"
# The previous "Periodo Mora" values (1802-1806) in the account-statement
# table are removed and replaced with the new periods, entered in the
# opposite order (1806 down to 1802) in the same cells (E16:E20).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "1806"
$ws.Range("E17").Value = "1805"
$ws.Range("E18").Value = "1804"
$ws.Range("E19").Value = "1803"
$ws.Range("E20").Value = "1802"
